$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) DevOps Internship at Provectus -> add "(Russia)"
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "DevOps Internship at Provectus 11/2021-02/2022",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "DevOps Internship at Provectus (Russia) 11/2021-02/2022", 2)

# ---------------------------------------------------------------------
# 2) New job entry: Machine Operator at EverLight (Taiwan), inserted
#    right after the "...alarm with SNS" bullet (end of the Provectus
#    internship section), before the horizontal rule / Technical
#    Competence heading.
# ---------------------------------------------------------------------
$findRng = $d.Content
$null = $findRng.Find.Execute("SNS")
$findRng.Collapse(0)
$insertPos = $findRng.Start
$insertRng = $d.Range($insertPos, $insertPos)

$jobXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Machine Operator at EverLight (Taiwan) 07/2015-07/2018</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">- Operate plating machines</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $insertRng.InsertXML($jobXml)

# ---------------------------------------------------------------------
# 3) Two new certification bullets under "Extra", inserted right after
#    the "Letter of recommendation from Yegor Bugayenko" bullet, before
#    "Human Languages:".
# ---------------------------------------------------------------------
$findRng2 = $d.Content
$null = $findRng2.Find.Execute("Letter of recommendation from Yegor Bugayenko")
$findRng2.Collapse(0)
$insertPos2 = $findRng2.Start
$insertRng2 = $d.Range($insertPos2, $insertPos2)

$certXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1009"/><w:ilvl w:val="0"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Computer Hardware I, II certification (Taiwan)</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1009"/><w:ilvl w:val="0"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Electronics I certification (Taiwan)</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $insertRng2.InsertXML($certXml)

Write-Output "edits applied"
